# Re-grades the quiz marksheet so that the student's float-typed answers are
# compared correctly against the correct answers (previously a float-vs-string
# comparison bug made the grader treat every question as "not attempted").
#
# After the fix:
#  - Group 1 (cols A/B, 25 questions) is graded: matching answers are filled
#    into the "Student Ans" column (A) and highlighted with the same green
#    "correct" style already used elsewhere on the sheet; unmatched ones stay
#    blank/black ("not attempted").
#  - Group 2 (cols D/E) only had 3 real questions; the unused D/E rows
#    (19-40) and the unused 3rd group (cols G/H) are cleared out entirely.
#  - The summary box (rows 10-12) is recomputed: Right/Not-Attempt/Max counts,
#    the marking scheme, the total score and the "scored/possible" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Summary box (rows 10-12)
# ---------------------------------------------------------------------

# Row 10 "No.": Right / Wrong / Not Attempt / Max counts
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 28

# Row 11 "Marking": points per right / wrong / not-attempted answer
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 "Total": total score and "scored/possible" label
$ws.Range("B12").Value = 76
$ws.Range("E12").Value = "76/112"

# A10/A11/A12 gain the same bold "mtitleStyle" formatting already used by the
# neighbouring label cells (copy format from A9, which already has it).
$labelStyleSrc = $ws.Range("A9")
foreach ($r in @(10, 11, 12)) {
    $labelStyleSrc.Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 2. Group 1 answers (cols A/B, rows 16-40): fill in "Student Ans" (A) for
#    every question the student got right, formatted with the same green
#    "correct" style used by B10 (correctStyle).
# ---------------------------------------------------------------------

$correctStyleSrc = $ws.Range("B10")
$rightAnswersCol = @{
    16 = "Option A"
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    21 = "Option C"
    22 = "Option D"
    23 = "Option D"
    25 = "Option A"
    26 = "Option C"
    27 = "Option A"
    30 = "Option B"
    32 = "Option C"
    33 = "Option D"
    36 = "Option A"
    37 = "Option A"
    39 = "Option D"
}

foreach ($r in $rightAnswersCol.Keys) {
    $correctStyleSrc.Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $rightAnswersCol[$r]
}

# ---------------------------------------------------------------------
# 3. Group 2 answers (cols D/E, rows 16-18): only 3 real questions, all
#    correct - fill "Student Ans" (D) to match "Correct Ans" (E), using the
#    same green "correct" style.
# ---------------------------------------------------------------------

$group2Answers = @{
    16 = "Option A"
    17 = "Option C"
    18 = "Option D"
}

foreach ($r in $group2Answers.Keys) {
    $correctStyleSrc.Copy() | Out-Null
    $ws.Cells.Item($r, 4).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 4).Value = $group2Answers[$r]
}

# ---------------------------------------------------------------------
# 4. Remove the now-unused cells: the rest of group 2 (D/E rows 19-40) and
#    the whole of the 3rd group (G/H rows 15-40).
# ---------------------------------------------------------------------

$ws.Range("D19:E40").Clear() | Out-Null
$ws.Range("G15:H40").Clear() | Out-Null
